# Weekly update: insert a new price observation row for
# "Feria Lagunitas de Puerto Montt - Haba" ahead of the existing row 79,
# shifting all subsequent rows (79-116 -> 80-117) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 79; Excel shifts rows 79:116 -> 80:117
# and extends the used range to row 117.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new data point.
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 44875
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112026
$ws.Range("G79").Value = "Haba"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 60
$ws.Range("K79").Value = 12000
$ws.Range("L79").Value = 12000
$ws.Range("M79").Value = 12000
$ws.Range("N79").Value = "$/saco 25 kilos"
$ws.Range("O79").Value = "Región del Maule"
$ws.Range("P79").Value = 480
$ws.Range("Q79").Value = 25
$ws.Range("R79").Value = "Hortaliza"
